# This workbook's data rows (2-18) are being reshuffled: the commit rotates
# the weekly records so that each row's D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), R (Origen), S (Precio $/Kg) and
# T (Kg / unidad) values move to a different row, per the mapping below
# (captured from a snapshot of the original data, then written back
# permuted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that participate in the row permutation.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# Snapshot the current (pre-edit) values for each affected column/row.
# (Value2 is used for reading because it reliably returns the raw
# number/string, whereas Value can return date/currency-wrapped variants.)
$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{}
    for ($r = 2; $r -le 18; $r++) {
        $orig[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# Mapping: new row -> source row (where the data used to live).
$mapping = @{
    2  = 14
    3  = 15
    4  = 16
    5  = 17
    6  = 2
    7  = 3
    8  = 18
    9  = 10
    10 = 11
    11 = 8
    12 = 9
    13 = 4
    14 = 5
    15 = 12
    16 = 13
    17 = 6
    18 = 7
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $orig[$col][$srcRow]
    }
}
